# Sharkspray_Testcases.xlsx - "Added code for hash"
# Fix a typo in row 3 (adhesiveType) and append a new test-case row (row 4)
# for the "oca2_p3" phase-3 (DMA + Compression + Tension) case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 3: adhesiveType had a truncated value missing the
# closing parenthesis; correct it to match row 2 / the new row 4.
$ws.Range("G3").Value = "PSA - Adhesive Transfer Tape (ATT)"

# --- Append new row 4 with the oca2_p3 test case.
# Fill order matters for shared-string allocation, so K4 is written
# before I4 to mirror the original authoring order.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "oca2_p3"
$ws.Range("C4").Value = "oca2_dma.xml"
$ws.Range("D4").Value = "oca2_compression.xml"
$ws.Range("E4").Value = "oca2_tension.xlsx"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "PSA - Adhesive Transfer Tape (ATT)"
$ws.Range("H4").Value = "Phase 3 (DMA + Compression + Tension)"
$ws.Range("K4").Value = "oca2_p3_ModelFiles.zip"
$ws.Range("I4").Value = "Tension"
$ws.Range("J4").Value = $true

# --- Column widths: widened (best-fit) to accommodate the new, longer
# content in columns D (compression_filename) and H (modelPhase).
$ws.Columns("D").ColumnWidth = 19.1
$ws.Columns("H").ColumnWidth = 33.25

# --- Selection moved off the edited cell, as recorded in the saved file.
$ws.Range("J16").Select()
